$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.858.98"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.634.34"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -1.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.78"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5020"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2556"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06343"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.30"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.235"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.865.66"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.630.98"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5388"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "0.0₅7861"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.14"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "25.887.60"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "196.24"
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.329"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.827"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.933"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.889"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.05"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1126"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.785"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.233"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04904"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.236"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.169"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.522"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.361"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8849"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.601"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("D38").Value = "1.140.21"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5513"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01563"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9999"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.664"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8078"
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.24"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "0.0₈119"
$ws.Range("E45").Value = "  +5.71%  "
$ws.Range("D46").Value = "1.777.13"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4504"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9990"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.28"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05051"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  -1.51%  "
